$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1952
$ws.Range("F5").Value = 452
$ws.Range("F6").Value = 1880
$ws.Range("F7").Value = 875
$ws.Range("F8").Value = 1381
$ws.Range("F9").Value = 8
$ws.Range("F10").Value = 1035
$ws.Range("F11").Value = 1035
$ws.Range("F13").Value = 3007
$ws.Range("F14").Value = 408
$ws.Range("F15").Value = 894
$ws.Range("F16").Value = 1166
$ws.Range("F17").Value = 606
$ws.Range("F20").Value = 1719
$ws.Range("F21").Value = 358
$ws.Range("F22").Value = 1290
$ws.Range("F23").Value = 228
$ws.Range("F24").Value = 608
$ws.Range("F25").Value = 508
$ws.Range("F26").Value = 1086
$ws.Range("F27").Value = 1577
$ws.Range("F28").Value = 1487
$ws.Range("F29").Value = 1353
$ws.Range("F30").Value = 389
$ws.Range("F31").Value = 1311
$ws.Range("F32").Value = 457
$ws.Range("F34").Value = 985
$ws.Range("F36").Value = 1872
$ws.Range("F37").Value = 494
$ws.Range("F38").Value = 1068
$ws.Range("F39").Value = 165
$ws.Range("F41").Value = 2317
$ws.Range("F42").Value = 161
$ws.Range("F43").Value = 900
$ws.Range("F44").Value = 2830
$ws.Range("F47").Value = 650

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 374
$ws.Range("F14").Value = 117376
$ws.Range("F22").Value = 294
$ws.Range("F23").Value = 5
$ws.Range("F28").Value = 62
$ws.Range("F29").Value = 73
$ws.Range("F30").Value = 73
$ws.Range("F34").Value = 157
$ws.Range("G34").Value = 280
$ws.Range("F41").Value = 164

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 3043
$ws.Range("F6").Value = 4861
$ws.Range("F7").Value = 185
$ws.Range("F9").Value = 697
$ws.Range("F10").Value = 965
$ws.Range("F11").Value = 560
$ws.Range("F12").Value = 652
$ws.Range("F13").Value = 1394
$ws.Range("F14").Value = 402
$ws.Range("F15").Value = 1324

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1952
$ws.Range("F4").Value = 697
$ws.Range("F5").Value = 965
$ws.Range("F6").Value = 560
$ws.Range("F8").Value = 1394
$ws.Range("F9").Value = 452
$ws.Range("F10").Value = 1880
$ws.Range("F11").Value = 875
$ws.Range("F12").Value = 1381
$ws.Range("F13").Value = 8
$ws.Range("F14").Value = 1035
$ws.Range("F15").Value = 1035
$ws.Range("F16").Value = 3007
$ws.Range("F18").Value = 408
$ws.Range("F19").Value = 894
$ws.Range("F20").Value = 1166
$ws.Range("F21").Value = 606
$ws.Range("F23").Value = 1719
$ws.Range("F25").Value = 358
$ws.Range("F26").Value = 374
$ws.Range("F27").Value = 608
$ws.Range("F28").Value = 508
$ws.Range("F29").Value = 1086
$ws.Range("F30").Value = 1577
$ws.Range("F31").Value = 1487
$ws.Range("F32").Value = 1353
$ws.Range("F33").Value = 389
$ws.Range("F35").Value = 1311
$ws.Range("F36").Value = 457
$ws.Range("F37").Value = 985
$ws.Range("F39").Value = 1872
$ws.Range("F40").Value = 73
$ws.Range("F41").Value = 1068
$ws.Range("F42").Value = 157
$ws.Range("G42").Value = 280
$ws.Range("F44").Value = 2317
$ws.Range("F45").Value = 161
$ws.Range("F46").Value = 900
$ws.Range("F47").Value = 2830
$ws.Range("F49").Value = 650
